$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; "E"=3; "G"=50.70817566666667; "H"=152.124527; "I"=0.5661129211027078; "J"=0.5661129211027077; "K"=3; "M"=153.5290173333333; "N"=460.587052; "O"=0.3172206968818489; "P"=0.317220696881849; "Q"=7785.176380869378; "R"=70066.5874278244; "S"=0.1795827353460201; "T"=0.1795827353460201}
    @{Row=3; "E"=3; "G"=50.70817566666667; "H"=152.124527; "I"=0.5661129211027078; "J"=0.5661129211027077; "K"=3; "M"=168.7997026666667; "N"=506.3991080000001; "O"=0.3487728915577651; "P"=0.3487728915577651; "Q"=8559.524975302436; "R"=77035.72477772193; "S"=0.1974448404412043; "T"=0.1974448404412043}
    @{Row=4; "E"=3; "G"=50.70817566666667; "H"=152.124527; "I"=0.5661129211027078; "J"=0.5661129211027077; "K"=3; "M"=68.09032333333333; "N"=204.27097; "O"=0.1406878008722904; "P"=0.1406878008722904; "Q"=3452.736076786799; "R"=31074.62469108119; "S"=0.07964518191532839; "T"=0.07964518191532839}
    @{Row=5; "E"=3; "G"=50.70817566666667; "H"=152.124527; "I"=0.5661129211027078; "J"=0.5661129211027077; "K"=3; "M"=93.562673; "N"=280.688019; "O"=0.1933186106880956; "P"=0.1933186106880956; "Q"=4744.392458326891; "R"=42699.53212494202; "S"=0.109440163400155; "T"=0.109440163400155}
    @{Row=6; "E"=3; "G"=17.08683666666667; "H"=51.26051; "I"=0.1907597520636141; "J"=0.1907597520636141; "K"=3; "M"=153.5290173333333; "N"=460.587052; "O"=0.3172206968818489; "P"=0.317220696881849; "Q"=2623.325242768502; "R"=23609.92718491652; "S"=0.06051294148662836; "T"=0.06051294148662838}
    @{Row=7; "E"=3; "G"=17.08683666666667; "H"=51.26051; "I"=0.1907597520636141; "J"=0.1907597520636141; "K"=3; "M"=168.7997026666667; "N"=506.3991080000001; "O"=0.3487728915577651; "P"=0.3487728915577651; "Q"=2884.252948847231; "R"=25958.27653962509; "S"=0.06653183032006903; "T"=0.06653183032006903}
    @{Row=8; "E"=3; "G"=17.08683666666667; "H"=51.26051; "I"=0.1907597520636141; "J"=0.1907597520636141; "K"=3; "M"=68.09032333333333; "N"=204.27097; "O"=0.1406878008722904; "P"=0.1406878008722904; "Q"=1163.448233377189; "R"=10471.0341003947; "S"=0.02683757001277322; "T"=0.02683757001277323}
    @{Row=9; "E"=3; "G"=17.08683666666667; "H"=51.26051; "I"=0.1907597520636141; "J"=0.1907597520636141; "K"=3; "M"=93.562673; "N"=280.688019; "O"=0.1933186106880956; "P"=0.1933186106880956; "Q"=1598.690111647743; "R"=14388.21100482969; "S"=0.03687741024414345; "T"=0.03687741024414345}
    @{Row=10; "E"=3; "G"=18.33915266666667; "H"=55.017458; "I"=0.2047407770084672; "J"=0.2047407770084672; "K"=3; "M"=153.5290173333333; "N"=460.587052; "O"=0.3172206968818489; "P"=0.317220696881849; "Q"=2815.592087639312; "R"=25340.32878875382; "S"=0.06494801196275718; "T"=0.06494801196275719}
    @{Row=11; "E"=3; "G"=18.33915266666667; "H"=55.017458; "I"=0.2047407770084672; "J"=0.2047407770084672; "K"=3; "M"=168.7997026666667; "N"=506.3991080000001; "O"=0.3487728915577651; "P"=0.3487728915577651; "Q"=3095.643517291941; "R"=27860.79165562747; "S"=0.07140803281702668; "T"=0.07140803281702668}
    @{Row=12; "E"=3; "G"=18.33915266666667; "H"=55.017458; "I"=0.2047407770084672; "J"=0.2047407770084672; "K"=3; "M"=68.09032333333333; "N"=204.27097; "O"=0.1406878008722904; "P"=0.1406878008722904; "Q"=1248.718834732696; "R"=11238.46951259426; "S"=0.02880452966620524; "T"=0.02880452966620524}
    @{Row=13; "E"=3; "G"=18.33915266666667; "H"=55.017458; "I"=0.2047407770084672; "J"=0.2047407770084672; "K"=3; "M"=93.562673; "N"=280.688019; "O"=0.1933186106880956; "P"=0.1933186106880956; "Q"=1715.860144048411; "R"=15442.7412964357; "S"=0.03958020256247807; "T"=0.03958020256247807}
    @{Row=14; "E"=3; "G"=3.438381; "H"=10.315143; "I"=0.03838654982521095; "J"=0.03838654982521095; "K"=3; "M"=153.5290173333333; "N"=460.587052; "O"=0.3172206968818489; "P"=0.317220696881849; "Q"=527.8912561476039; "R"=4751.021305328436; "S"=0.01217700808644323; "T"=0.01217700808644324}
    @{Row=15; "E"=3; "G"=3.438381; "H"=10.315143; "I"=0.03838654982521095; "J"=0.03838654982521095; "K"=3; "M"=168.7997026666667; "N"=506.3991080000001; "O"=0.3487728915577651; "P"=0.3487728915577651; "Q"=580.3976904547161; "R"=5223.579214092445; "S"=0.01338818797946505; "T"=0.01338818797946505}
    @{Row=16; "E"=3; "G"=3.438381; "H"=10.315143; "I"=0.03838654982521095; "J"=0.03838654982521095; "K"=3; "M"=68.09032333333333; "N"=204.27097; "O"=0.1406878008722904; "P"=0.1406878008722904; "Q"=234.12047403319; "R"=2107.08426629871; "S"=0.005400519277983531; "T"=0.005400519277983533}
    @{Row=17; "E"=3; "G"=3.438381; "H"=10.315143; "I"=0.03838654982521095; "J"=0.03838654982521095; "K"=3; "M"=93.562673; "N"=280.688019; "O"=0.1933186106880956; "P"=0.1933186106880956; "Q"=321.704117152413; "R"=2895.337054371717; "S"=0.007420834481319142; "T"=0.007420834481319142}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Range("E$row").Value = $r["E"]
    $ws.Range("G$row").Value = $r["G"]
    $ws.Range("H$row").Value = $r["H"]
    $ws.Range("I$row").Value = $r["I"]
    $ws.Range("J$row").Value = $r["J"]
    $ws.Range("K$row").Value = $r["K"]
    $ws.Range("M$row").Value = $r["M"]
    $ws.Range("N$row").Value = $r["N"]
    $ws.Range("O$row").Value = $r["O"]
    $ws.Range("P$row").Value = $r["P"]
    $ws.Range("Q$row").Value = $r["Q"]
    $ws.Range("R$row").Value = $r["R"]
    $ws.Range("S$row").Value = $r["S"]
    $ws.Range("T$row").Value = $r["T"]
}
